# Actualización desde MV -datos-
# Adds a new quarterly data row (01-07-2021) to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 68

# Column A holds period labels stored as text (e.g. "01-04-2021"); force
# text format first so Excel doesn't auto-convert the dd-mm-yyyy-looking
# string into a date serial, then drop the format again so the cell keeps
# the workbook's default (unstyled) look, matching the other data rows.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "01-07-2021"
$ws.Cells.Item($newRow, 1).ClearFormats()

$values = @{
    2  = 118912
    3  = 205751
    4  = 459691
    5  = 253940
    6  = -509018
    7  = 1502936
    8  = 2011954
    9  = 28689
    10 = 393490
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 8029598
    17 = 5566866
    18 = 5928008
    19 = 5925987
    20 = 2021
    21 = 361142
    22 = 2555672
    23 = 4650816
    24 = 4650816
    25 = 0
    26 = 2095144
    27 = -92940
    28 = -7910687
}

foreach ($col in $values.Keys) {
    $ws.Cells.Item($newRow, $col).Value = $values[$col]
}

$wb.Save()
